$wb = $excel.ActiveWorkbook

# --- Core data edit (QA environment bank details for paymentMethod sheet) ---
$paymentMethod = $wb.Worksheets.Item("paymentMethod")
$paymentMethod.Range("C6").Value = "ANDHRA BANK-Andhra Bank RTC Busstand"
$paymentMethod.Range("D6").Value = "110710011005899"

# Selection moved from E7 to E6 on the paymentMethod sheet
[void]$paymentMethod.Range("E6").Select()

# --- Minor column-width relayout that accompanied the edit (all sheets) ---
$chequeDetails = $wb.Worksheets.Item("chequeDetails")
$chequeDetails.Range("A1").EntireColumn.ColumnWidth = 24.5
$chequeDetails.Range("B1").EntireColumn.ColumnWidth = 17.6666666666667
$chequeDetails.Range("C1").EntireColumn.ColumnWidth = 9.16666666666667
$chequeDetails.Range("D1").EntireColumn.ColumnWidth = 27.8333333333333

$challanHeaderDetails = $wb.Worksheets.Item("challanHeaderDetails")
$challanHeaderDetails.Range("A1").EntireColumn.ColumnWidth = 9.83333333333333
$challanHeaderDetails.Range("B1").EntireColumn.ColumnWidth = 20.8333333333333
$challanHeaderDetails.Range("C1").EntireColumn.ColumnWidth = 21.3333333333333
$challanHeaderDetails.Range("D1").EntireColumn.ColumnWidth = 20.6666666666667
$challanHeaderDetails.Range("E1").EntireColumn.ColumnWidth = 20.6666666666667
$challanHeaderDetails.Range("F1").EntireColumn.ColumnWidth = 21.3333333333333
$challanHeaderDetails.Range("G1").EntireColumn.ColumnWidth = 21.3333333333333
$challanHeaderDetails.Range("H1").EntireColumn.ColumnWidth = 9.83333333333333

$approvalDetails = $wb.Worksheets.Item("approvalDetails")
$approvalDetails.Range("A1").EntireColumn.ColumnWidth = 21
$approvalDetails.Range("B1").EntireColumn.ColumnWidth = 32.3333333333333
$approvalDetails.Range("C1").EntireColumn.ColumnWidth = 32.5
$approvalDetails.Range("D1").EntireColumn.ColumnWidth = 51.3333333333333

$paymentMethod.Range("A1").EntireColumn.ColumnWidth = 16.5
$paymentMethod.Range("B1").EntireColumn.ColumnWidth = 14
$paymentMethod.Range("C1").EntireColumn.ColumnWidth = 14.5
$paymentMethod.Range("D1").EntireColumn.ColumnWidth = 15.6666666666667
